$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update header row (row 1) question text ---
$ws1.Range("B1").Value = "Cuantos tachos de basura hay en su cuadra?"
$ws1.Range("C1").Value = "Cuantos arboles hay en su cuadra? (aprox)"

# --- Build the new response rows (4, 5, 6) by first copying the ---
# --- formatting from an existing response row, while both row 2 and ---
# --- row 3 still hold their original formatting. ---
$ws1.Range("A3:E3").Copy($ws1.Range("A4:E4"))
$ws1.Range("A3:E3").Copy($ws1.Range("A5:E5"))
$ws1.Range("A3:E3").Copy($ws1.Range("A6:E6"))

# --- Now overwrite the new rows with their actual values. ---
$ws1.Range("A4").Value = 45248.92231679398
$ws1.Range("B4").Value = 4
$ws1.Range("C4").Value = 7
$ws1.Range("D4").Value = "No"
$ws1.Range("E4").Value = "eitanluc@gmail.com"

$ws1.Range("A5").Value = 45248.922446145836
$ws1.Range("B5").Value = 3
$ws1.Range("C5").Value = 10
$ws1.Range("D5").Value = "No"
$ws1.Range("E5").Value = "eitanbaserow@gmail.com"

$ws1.Range("A6").Value = 45248.92258255787
$ws1.Range("B6").Value = 4
$ws1.Range("C6").Value = 7
$ws1.Range("D6").Value = "Si, pero en mal estado"
$ws1.Range("E6").Value = "eitanluc@gmail.com"

# --- Clear out the old response rows 2 and 3: drop the old B:E cells ---
# --- entirely, and clear the A cell's value while keeping its style. ---
$ws1.Range("B2:E2").Clear()
$ws1.Range("A2").ClearContents()

$ws1.Range("B3:E3").Clear()
$ws1.Range("A3").ClearContents()
